$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'327.81"
$ws.Range('E2').Value = "'1.07%"
$ws.Range('G2').Value = "'6"
$ws.Range('B3').Value = "'HuobiToken"
$ws.Range('C3').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D3').Value = "'5.506"
$ws.Range('E3').Value = "'0.12%"
$ws.Range('G3').Value = "'6"
$ws.Range('B4').Value = "'Cronos"
$ws.Range('C4').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D4').Value = "'0.08014"
$ws.Range('E4').Value = "'-0.24%"
$ws.Range('G4').Value = "'6"
$ws.Range('B5').Value = "'FTXToken"
$ws.Range('C5').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D5').Value = "'2.018"
$ws.Range('E5').Value = "'7.10%"
$ws.Range('G5').Value = "'6"
$ws.Range('B6').Value = "'BTSEToken"
$ws.Range('C6').Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range('D6').Value = "'2.569"
$ws.Range('E6').Value = "'-3.48%"
$ws.Range('G6').Value = "'6"
$ws.Range('B7').Value = "'MXToken"
$ws.Range('C7').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D7').Value = "'0.9489"
$ws.Range('E7').Value = "'0.97%"
$ws.Range('G7').Value = "'6"
$ws.Range('B8').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C8').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D8').Value = "'0.1122"
$ws.Range('E8').Value = "'-4.76%"
$ws.Range('G8').Value = "'6"
$ws.Range('B9').Value = "'WazirX"
$ws.Range('C9').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D9').Value = "'0.1859"
$ws.Range('E9').Value = "'-1.08%"
$ws.Range('G9').Value = "'6"
$ws.Range('B10').Value = "'MCDex"
$ws.Range('C10').Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range('D10').Value = "'10.61"
$ws.Range('E10').Value = "'25.10%"
$ws.Range('G10').Value = "'6"
$ws.Range('B11').Value = "'MandalaExchangeToken"
$ws.Range('C11').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D11').Value = "'0.09921"
$ws.Range('E11').Value = "'-0.38%"
$ws.Range('G11').Value = "'6"
$ws.Range('B12').Value = "'BitrueCoin"
$ws.Range('C12').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D12').Value = "'0.04633"
$ws.Range('E12').Value = "'11.43%"
$ws.Range('G12').Value = "'6"
$ws.Range('B13').Value = "'BitMartToken"
$ws.Range('C13').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D13').Value = "'0.1064"
$ws.Range('E13').Value = "'-0.15%"
$ws.Range('G13').Value = "'6"
$ws.Range('B14').Value = "'BitForexToken"
$ws.Range('C14').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D14').Value = "'0.001263"
$ws.Range('E14').Value = "'-0.82%"
$ws.Range('G14').Value = "'6"
$ws.Range('B15').Value = "'CoinExToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D15').Value = "'0.04076"
$ws.Range('E15').Value = "'-4.26%"
$ws.Range('G15').Value = "'6"
$ws.Range('B16').Value = "'TigerCash"
$ws.Range('C16').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('D16').Value = "'0.005941"
$ws.Range('E16').Value = "'-0.82%"
$ws.Range('G16').Value = "'6"
$ws.Range('B17').Value = "'OKB"
$ws.Range('C17').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D17').Value = "'43.81"
$ws.Range('E17').Value = "'-1.80%"
$ws.Range('G17').Value = "'6"
$ws.Range('B18').Value = "'LEO"
$ws.Range('C18').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('D18').Value = "'3.356"
$ws.Range('E18').Value = "'-6.71%"
$ws.Range('G18').Value = "'6"
$ws.Range('B19').Value = "'GateToken"
$ws.Range('C19').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D19').Value = "'4.319"
$ws.Range('E19').Value = "'0.30%"
$ws.Range('G19').Value = "'6"
$ws.Range('D20').Value = "'0.3476"
$ws.Range('E20').Value = "'-0.30%"
$ws.Range('G20').Value = "'6"
$ws.Range('D21').Value = "'0.1406"
$ws.Range('E21').Value = "'2.37%"
$ws.Range('G21').Value = "'6"
$ws.Range('D22').Value = "'0.2545"
$ws.Range('E22').Value = "'-3.75%"
$ws.Range('G22').Value = "'6"
$ws.Range('E23').Value = "'1.34%"
$ws.Range('G23').Value = "'6"
$ws.Range('D24').Value = "'0.004337"
$ws.Range('E24').Value = "'-2.61%"
$ws.Range('G24').Value = "'6"
$ws.Range('E25').Value = "'-3.84%"
$ws.Range('G25').Value = "'6"
$ws.Range('D26').Value = "'0.0003742"
$ws.Range('E26').Value = "'-6.54%"
$ws.Range('G26').Value = "'6"
$ws.Range('G27').Value = "'6"
$ws.Range('G28').Value = "'6"
$ws.Range('G29').Value = "'6"
$ws.Range('G30').Value = "'6"
$ws.Range('G31').Value = "'6"
$ws.Range('G32').Value = "'6"
$ws.Range('G33').Value = "'6"
$ws.Range('G34').Value = "'6"
$ws.Range('G35').Value = "'6"
$ws.Range('G36').Value = "'6"
$ws.Range('G37').Value = "'6"
$ws.Range('D38').Value = "'0.02577"
$ws.Range('E38').Value = "'-1.77%"
$ws.Range('G38').Value = "'6"
$ws.Range('D39').Value = "'0.05631"
$ws.Range('E39').Value = "'2.84%"
$ws.Range('G39').Value = "'6"
$ws.Range('D40').Value = "'0.007534"
$ws.Range('E40').Value = "'-1.92%"
$ws.Range('G40').Value = "'6"
$ws.Range('E41').Value = "'0.29%"
$ws.Range('G41').Value = "'6"
$ws.Range('D42').Value = "'0.007516"
$ws.Range('E42').Value = "'4.33%"
$ws.Range('G42').Value = "'6"
$ws.Range('D43').Value = "'0.002014"
$ws.Range('E43').Value = "'-2.02%"
$ws.Range('G43').Value = "'6"
$ws.Range('D44').Value = "'0.008381"
$ws.Range('E44').Value = "'-8.80%"
$ws.Range('G44').Value = "'6"
$ws.Range('D45').Value = "'0.00007097"
$ws.Range('E45').Value = "'-0.33%"
$ws.Range('G45').Value = "'6"
$ws.Range('D46').Value = "'0.00000000750"
$ws.Range('E46').Value = "'-0.42%"
$ws.Range('G46').Value = "'6"
$ws.Range('E47').Value = "'54.88%"
$ws.Range('G47').Value = "'6"
$ws.Range('D48').Value = "'0.003454"
$ws.Range('E48').Value = "'-0.99%"
$ws.Range('G48').Value = "'6"
$ws.Range('D49').Value = "'0.00002099"
$ws.Range('E49').Value = "'-0.42%"
$ws.Range('G49').Value = "'6"
$ws.Range('D50').Value = "'0.0001999"
$ws.Range('E50').Value = "'-0.42%"
$ws.Range('G50').Value = "'6"
$ws.Range('G51').Value = "'6"
